$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Update the surviving sheet's data (was Sheet3) before removing Sheet1 ---
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A2").Value = "standard_user"
$ws3.Range("A3").Value = "problem_user"
$ws3.Range("A4").Value = "error_user"

# --- Remove Sheet1 entirely ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Delete()

# --- Re-fetch Sheet3 (now the only / active sheet) and update its selection ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A5").Select()
